# Falcons Team Data.xlsx — "Logged Week 15 and simulated Week 16"
# Appends per-play logged values for weeks 15/16 to the long running
# number-list cells on YDS and ST, and bumps the season summary totals
# on OFF, DEF, ST, TURNS and PEN accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet: running per-play yardage logs (Rush/Pass, OFF/DEF)
# ---------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 4 0 0 7 2 7 1 3 2 3 3 6 0 0 2 -1 5 4 8 2 1 3"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 7 7 14 5 0 13 20 10 21 49 2 4 49 5 14 3 4 5 4"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 2 5 5 3 2 2 17 2 3 6 0 3 10 0 5 6 12 8 13 13 1 3 5 1 1 16 3 3 2 6 2 2"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 6 25 6 28 3 36 9 9 30 3 15 13 21 15 3 4 8 1"

# ---------------------------------------------------------------
# OFF sheet: season totals (row 2 = Home, row 3 = Road)
# ---------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 379
$offWs.Range("D2").Value = 16
$offWs.Range("F2").Value = 92
$offWs.Range("G2").Value = 84
$offWs.Range("I2").Value = 14
$offWs.Range("J2").Value = 50
$offWs.Range("N2").Value = 37
$offWs.Range("O2").Value = 45

$offWs.Range("B3").Value = 13
$offWs.Range("C3").Value = 303
$offWs.Range("D3").Value = 9
$offWs.Range("E3").Value = 58
$offWs.Range("F3").Value = 233
$offWs.Range("G3").Value = 60
$offWs.Range("I3").Value = 127
$offWs.Range("J3").Value = 109
$offWs.Range("L3").Value = 549
$offWs.Range("M3").Value = 369
$offWs.Range("Q3").Value = 986

# ---------------------------------------------------------------
# DEF sheet: season totals (row 2 = Home, row 3 = Road)
# ---------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 353
$defWs.Range("D2").Value = 22
$defWs.Range("F2").Value = 93
$defWs.Range("G2").Value = 114
$defWs.Range("J2").Value = 51

$defWs.Range("B3").Value = 19
$defWs.Range("C3").Value = 366
$defWs.Range("E3").Value = 65
$defWs.Range("F3").Value = 207
$defWs.Range("G3").Value = 68
$defWs.Range("I3").Value = 105
$defWs.Range("J3").Value = 106
$defWs.Range("L3").Value = 588
$defWs.Range("M3").Value = 402
$defWs.Range("Q3").Value = 1017

# ---------------------------------------------------------------
# ST sheet: season totals + running KO/PT distance logs
# ---------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 144
$stWs.Range("D2").Value = 112
$stWs.Range("F2").Value = 86
$stWs.Range("G2").Value = 82
$stWs.Range("J2").Value = 57
$stWs.Range("K2").Value = 55
$stWs.Range("L2").Value = 27
$stWs.Range("M2").Value = 20

$stWs.Range("B3").Value = 76

$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 64"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 53 46"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 10 -2"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 15"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0 0 23 0"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 20 26 33 24"

# ---------------------------------------------------------------
# TURNS sheet: Road FMBLg total
# ---------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("E3").Value = 14

# ---------------------------------------------------------------
# PEN sheet: False start count + Intentional grounding yards
# ---------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value = 39
$penWs.Range("D4").Value = 21
